$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-21 19:48:38"
$ws.Range("E3").Value = "2026-02-21 19:48:41"
$ws.Range("E4").Value = "2026-02-21 19:48:44"
$ws.Range("H4").Value = "'71%"
$ws.Range("O4").Value = "9.7 °C"
$ws.Range("E5").Value = "2026-02-21 19:48:46"
$ws.Range("E6").Value = "2026-02-21 19:48:49"
$ws.Range("E7").Value = "2026-02-21 19:48:52"
$ws.Range("H7").Value = "'55%"
$ws.Range("E8").Value = "2026-02-21 19:48:54"
$ws.Range("O8").Value = "11.0 °C"
$ws.Range("E9").Value = "2026-02-21 19:48:57"
$ws.Range("N9").Value = "7.8 °C 19:29 TU"
$ws.Range("O9").Value = "13.8 °C"
$ws.Range("E10").Value = "2026-02-21 19:49:00"
$ws.Range("O10").Value = "8.7 °C"
$ws.Range("E11").Value = "2026-02-21 19:49:02"
$ws.Range("O11").Value = "9.2 °C"
$ws.Range("E12").Value = "2026-02-21 19:49:05"
$ws.Range("H12").Value = "'60%"
$ws.Range("E13").Value = "2026-02-21 19:49:07"
$ws.Range("H13").Value = "'60%"
$ws.Range("E14").Value = "2026-02-21 19:49:10"
$ws.Range("H14").Value = "'68%"
$ws.Range("N14").Value = "5.9 °C 19:28 TU"
$ws.Range("O14").Value = "11.6 °C"
$ws.Range("E15").Value = "2026-02-21 19:49:13"
$ws.Range("H15").Value = "'53%"
$ws.Range("O15").Value = "13.5 °C"
$ws.Range("E16").Value = "2026-02-21 19:49:15"
$ws.Range("E17").Value = "2026-02-21 19:49:18"
$ws.Range("E18").Value = "2026-02-21 19:49:20"
$ws.Range("E19").Value = "2026-02-21 19:49:23"
$ws.Range("E20").Value = "2026-02-21 19:49:26"
$ws.Range("E21").Value = "2026-02-21 19:49:28"
$ws.Range("O21").Value = "7.7 °C"
$ws.Range("E22").Value = "2026-02-21 19:49:31"
$ws.Range("O22").Value = "1.9 °C"
$ws.Range("E23").Value = "2026-02-21 19:49:34"
$ws.Range("E24").Value = "2026-02-21 19:49:36"
$ws.Range("E25").Value = "2026-02-21 19:49:39"
$ws.Range("E26").Value = "2026-02-21 19:49:42"
$ws.Range("G26").Value = "3 cm"
$ws.Range("H26").Value = "'36%"
$ws.Range("J26").Value = "1027.2 hPa"
$ws.Range("K26").Value = "15.5 MJ/m2"
$ws.Range("M26").Value = "14.5 °C 11:40 TU"
$ws.Range("O26").Value = "9.8 °C"
$ws.Range("E27").Value = "2026-02-21 19:49:45"
$ws.Range("E28").Value = "2026-02-21 19:49:47"
$ws.Range("E29").Value = "2026-02-21 19:49:50"
$ws.Range("H29").Value = "'65%"
$ws.Range("O29").Value = "11.9 °C"
$ws.Range("E30").Value = "2026-02-21 19:49:53"
$ws.Range("H30").Value = "'65%"
$ws.Range("O30").Value = "11.8 °C"
$ws.Range("E31").Value = "2026-02-21 19:49:55"
$ws.Range("E32").Value = "2026-02-21 19:49:58"
$ws.Range("O32").Value = "5.6 °C"
$ws.Range("E33").Value = "2026-02-21 19:50:01"
$ws.Range("J33").Value = "1030.2 hPa"
$ws.Range("O33").Value = "6.7 °C"
$ws.Range("E34").Value = "2026-02-21 19:50:04"
$ws.Range("H34").Value = "'38%"
$ws.Range("N34").Value = "0.0 °C 19:25 TU"
$ws.Range("O34").Value = "4.8 °C"
$ws.Range("E35").Value = "2026-02-21 19:50:06"
$ws.Range("O35").Value = "7.8 °C"
$ws.Range("E36").Value = "2026-02-21 19:50:09"
$ws.Range("H36").Value = "'56%"
$ws.Range("O36").Value = "13.5 °C"
$ws.Range("E37").Value = "2026-02-21 19:50:11"
$ws.Range("H37").Value = "'73%"
$ws.Range("E38").Value = "2026-02-21 19:50:14"
$ws.Range("E39").Value = "2026-02-21 19:50:17"
$ws.Range("E40").Value = "2026-02-21 19:50:19"
$ws.Range("H40").Value = "'51%"
$ws.Range("E41").Value = "2026-02-21 19:50:22"
$ws.Range("H41").Value = "'66%"
$ws.Range("E42").Value = "2026-02-21 19:50:24"
$ws.Range("H42").Value = "'73%"
$ws.Range("O42").Value = "11.0 °C"
$ws.Range("E43").Value = "2026-02-21 19:50:27"
$ws.Range("E44").Value = "2026-02-21 19:50:29"
$ws.Range("H44").Value = "'39%"
$ws.Range("N44").Value = "-0.4 °C 19:13 TU"
$ws.Range("O44").Value = "2.4 °C"
$ws.Range("E45").Value = "2026-02-21 19:50:32"
$ws.Range("E46").Value = "2026-02-21 19:50:35"
$ws.Range("H46").Value = "'68%"
